# dev update -- almost ready to merge back to the head
#
# Adds the second batch of "line-in" trial results (rows 31-41), the
# summary stats (AVERAGE/STDEV) under them, a new note, and a blank
# template block (rows 46-59 + 62) ready for the next round of trials.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Notes column (N) for the existing trial block (rows 29-34)
# ---------------------------------------------------------------------
$ws.Cells.Item(30,14).Value = "'-using line in with minimum volume to trigger collection, 3 dB gain on line input"
$ws.Cells.Item(31,14).Value = "'-direct speaker connection with out using microphone and pre-amp"
$ws.Cells.Item(32,14).Value = "'-comparison fft with microphone signal in soi[n]"
$ws.Cells.Item(33,14).Value = "'-tests run with sounds on loop on vlc media player--pseudo random"
$ws.Cells.Item(34,14).Value = "'-using new reset algorithm to avoid stopping and starting dsk"

# ---------------------------------------------------------------------
# Fill in the D/G/J date columns' number format up front (rows 31-41)
# so new cells created below pick up the existing "d-mmm" style used
# throughout the sheet instead of inventing a new one.
# ---------------------------------------------------------------------
$ws.Range("G31:G41").NumberFormat = "d-mmm"
$ws.Range("J31:J41").NumberFormat = "d-mmm"

# ---------------------------------------------------------------------
# Row-by-row trial data for the four mic/line-in trial groups
# (columns C/F/I/L = match %, columns D/G/J = their date stamps)
# ---------------------------------------------------------------------
$fVals = @(53,58,58,59,52,58,59,50,32,27,60)
$iVals = @(8,29,11,12,38,40,24,12,39,9,12)
$lVals = @(9,13,13,14,16,15,17,18,9,15,11)

for ($i = 0; $i -lt $fVals.Length; $i++) {
    $r = 31 + $i
    $ws.Cells.Item($r,6).Value  = $fVals[$i]
    $ws.Cells.Item($r,7).Value  = 41005
    $ws.Cells.Item($r,9).Value  = $iVals[$i]
    $ws.Cells.Item($r,10).Value = 41005
    $ws.Cells.Item($r,12).Value = $lVals[$i]
}

# Row 41 also gains its own A/C entry (date + match %), like rows 31-40
$ws.Cells.Item(41,1).Value = 41005
$ws.Cells.Item(41,3).Value = 62

# ---------------------------------------------------------------------
# Summary rows: averages, std-devs, and a closing note
# ---------------------------------------------------------------------
$ws.Cells.Item(42,3).Formula  = "=AVERAGE(C31:C41)"
$ws.Cells.Item(42,6).Formula  = "=AVERAGE(F31:F41)"
$ws.Cells.Item(42,9).Formula  = "=AVERAGE(I31:I41)"
$ws.Cells.Item(42,12).Formula = "=AVERAGE(L31:L41)"

$ws.Cells.Item(43,3).Formula  = "=STDEV(C30:C41)"
$ws.Cells.Item(43,6).Formula  = "=STDEV(F30:F41)"
$ws.Cells.Item(43,9).Formula  = "=STDEV(I30:I41)"
$ws.Cells.Item(43,12).Formula = "=STDEV(L30:L41)"

$ws.Cells.Item(44,3).Value = "'-set decision at 55 and you're wrong 16% of the time"

# ---------------------------------------------------------------------
# Next trial block's header rows (copy of the row 28/29 & 11/12 layout)
# ---------------------------------------------------------------------
$ws.Cells.Item(46,1).Value = "deep_growl_converted.wav"
$ws.Cells.Item(46,4).Value = "Child1.wav"
$ws.Cells.Item(46,7).Value = "Dogsbark.wav"
$ws.Cells.Item(46,10).Value = "THEFORCE.wav"

$headerCols = @(1,4,7,10)
foreach ($c in $headerCols) {
    $ws.Cells.Item(47,$c).Value   = "date"
    $ws.Cells.Item(47,$c+1).Value = "test #"
    $ws.Cells.Item(47,$c+2).Value = "match %"
}

# Blank, pre-formatted rows awaiting the next round of trials.
# Rows 48-53 keep all four date columns (A/D/G/J); rows 54-59 drop the
# "Child1.wav" column (D), matching the trimmed layout further down.
$ws.Range("A48:A53").NumberFormat = "d-mmm"
$ws.Range("D48:D53").NumberFormat = "d-mmm"
$ws.Range("G48:G59").NumberFormat = "d-mmm"
$ws.Range("J48:J59").NumberFormat = "d-mmm"
$ws.Range("A48:A59").NumberFormat = "d-mmm"

# Empty quote-prefixed cell under the closing note, ready for more text.
$ws.Cells.Item(62,3).Value = "'x"
$ws.Cells.Item(62,3).ClearContents()

# ---------------------------------------------------------------------
# View state: scrolled down to the new block, C48 selected
# ---------------------------------------------------------------------
$ws.Range("C48").Select()
